# Actualizacion automatica - renombrado de medidas "porcentaje-de-poblacion-*"
# a "de-poblacion-*" y conversion de la columna "Tasa de feminidad" (AA) de
# medida a dimension, con su correspondiente fichero de mapeo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 2: identificadores "slug" de las medidas de porcentaje ---
$ws.Range("C2").Value  = "de-poblacion-menor-de-25"
$ws.Range("I2").Value  = "de-poblacion-de-65-y-mas-anos"
$ws.Range("Q2").Value  = "de-poblacion-de-20-a-64-anos"
$ws.Range("R2").Value  = "de-poblacion-menor-de-45"
$ws.Range("X2").Value  = "de-poblacion-menor-de-15"
$ws.Range("Z2").Value  = "de-poblacion-de-0-a-19-anos"
$ws.Range("AF2").Value = "de-poblacion-menor-de-35"

# --- Fila 3: identificadores "iaest-measure:" de las mismas medidas ---
$ws.Range("C3").Value  = "iaest-measure:de-poblacion-menor-de-25"
$ws.Range("I3").Value  = "iaest-measure:de-poblacion-de-65-y-mas-anos"
$ws.Range("Q3").Value  = "iaest-measure:de-poblacion-de-20-a-64-anos"
$ws.Range("R3").Value  = "iaest-measure:de-poblacion-menor-de-45"
$ws.Range("X3").Value  = "iaest-measure:de-poblacion-menor-de-15"
$ws.Range("Z3").Value  = "iaest-measure:de-poblacion-de-0-a-19-anos"
$ws.Range("AF3").Value = "iaest-measure:de-poblacion-menor-de-35"

# --- Columna AA ("Tasa de feminidad"): pasa de medida a dimension ---
$ws.Range("AA3").Value = "iaest-dimension:tasa-de-feminidad"
$ws.Range("AA4").Value = "dim"
$ws.Range("AA5").Value = "skos:Concept"
$ws.Range("AA6").Value = "mapping-tasa-de-feminidad.xlsx"
